# "Adding contact details class"
# Adds a customerName column to the Login sheet and several new contact-detail
# columns (OtherId, DriverLicenseNo., SNNnumber, SINnumber, NickName,
# Militreyservice, DriverLicenseDate) to the createCustomer sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# createCustomer sheet - new header row (H1:N1) and data row (H2:N2)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("createCustomer")

$ws3.Range("H1").Value = "OtherId"
$ws3.Range("I1").Value = "DriverLicenseNo."
$ws3.Range("J1").Value = "SNNnumber"
$ws3.Range("K1").Value = "SINnumber"
$ws3.Range("L1").Value = "NickName"
$ws3.Range("M1").Value = "Militreyservice"

$ws3.Range("H2").Value = "1d0123"
$ws3.Range("I2").Value = "dl12345"
$ws3.Range("L2").Value = "Nick"
$ws3.Range("M2").Value = "None"

# Entered with a leading apostrophe so Excel keeps them as text (quote-prefix
# style) instead of converting to numbers.
$ws3.Range("J2").Value = "'1234"
$ws3.Range("K2").Value = "'5678"

# ---------------------------------------------------------------------------
# Login sheet - new customerName column (C1:C2)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Login")
$ws1.Range("C1").Value = "customerName"
$ws1.Range("C2").Value = "Peter"

# ---------------------------------------------------------------------------
# createCustomer sheet - DriverLicenseDate column (N1:N2)
# ---------------------------------------------------------------------------
$ws3.Range("N1").Value = "DriverLicenseDate"
# Also entered with a leading apostrophe (text, quote-prefixed) but formatted
# with the built-in short-date number format.
$ws3.Range("N2").Value = "'01/01/2000"
$ws3.Range("N2").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# Column widths (best effort - closest achievable autosize of the new columns)
# ---------------------------------------------------------------------------
$ws1.Columns.Item(3).ColumnWidth = 13.7576

$ws3.Columns.Item(8).ColumnWidth = 6.9221
$ws3.Columns.Item(9).ColumnWidth = 15.4218
$ws3.Columns.Item(10).ColumnWidth = 10.9221
$ws3.Columns.Item(11).ColumnWidth = 10.0867
$ws3.Columns.Item(12).ColumnWidth = 9.2573
$ws3.Columns.Item(13).ColumnWidth = 13.5924
$ws3.Columns.Item(14).ColumnWidth = 15.4218

# ---------------------------------------------------------------------------
# Selections
# ---------------------------------------------------------------------------
$ws1.Range("D2").Select() | Out-Null
$ws3.Range("J10").Select() | Out-Null
